$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.261.98'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '2.604.66'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''583.58'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = '''142.85'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("D9").Value = '''6.50'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -2.25%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").Value = '''0.372'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '3.064.22'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '''24.60'
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").Value = '60.262.03'
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").Value = '''0.0000141'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '2.608.80'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '''11.35'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("D20").Value = '''347.38'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").Value = '''6.90'
$ws.Range("E21").Value = '  -2.46%  '
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '''0.536'
$ws.Range("E23").Value = '  +3.87%  '
$ws.Range("D24").Value = '''63.79'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("E27").Value = '  +2.56%  '
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("D30").Value = '''169.19'
$ws.Range("E30").Value = '  +4.61%  '
$ws.Range("D31").Value = '''6.39'
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '''19.42'
$ws.Range("D34").Value = '''1.32'
$ws.Range("E34").Value = '  +10.16%  '
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").Value = '''0.992'
$ws.Range("E36").Value = '  +3.29%  '
$ws.Range("E37").Value = '  +2.47%  '
$ws.Range("D38").Value = '''317.19'
$ws.Range("E38").Value = '  +6.45%  '
$ws.Range("D39").Value = '''38.24'
$ws.Range("E39").Value = '  +1.43%  '
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("D41").Value = '''0.850'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").Value = '''135.70'
$ws.Range("E42").Value = '  -2.55%  '
$ws.Range("D43").Value = '''0.0994'
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").Value = '''19.95'
$ws.Range("E45").Value = '  +1.30%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").Value = '''19.96'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").Value = '''10.73'
$ws.Range("E51").Value = '  +0.26%  '
